# Update testcase live event and sticky ads
# Insert a new worksheet "List_Tab_Live_Event" right after "Continue_Watching"
# (becomes the 3rd sheet / new active tab) and populate it with the
# Index/Tabs table used by the other "List_Tab_*" style sheets.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("Continue_Watching")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "List_Tab_Live_Event"

# --- Data -------------------------------------------------------------
$newSheet.Range("A1").Value = "Index"
$newSheet.Range("B1").Value = "Tabs"

$newSheet.Range("A2").Value = "'0"
$newSheet.Range("B2").Value = "Live Event"

$newSheet.Range("A3").Value = "'1"
$newSheet.Range("B3").Value = "Missed Event"

# --- Formatting ---------------------------------------------------------
# Re-use the existing bold "header" style (fontId bold + themed fill)
# already present in the workbook by copying it from another sheet that
# carries the same Index/Value header look, instead of re-building it
# (which would create duplicate style entries).
$styleSource = $wb.Worksheets.Item("Homepage_Menu_Bawah")
$styleSource.Range("A1").Copy() | Out-Null
$newSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column B sizing similar to the other list sheets (matches the
# auto-fit width Excel computes for this column's longest entry,
# "Missed Event")
$newSheet.Columns.Item(2).ColumnWidth = 11.81640625

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

$newSheet.Range("A1").Select() | Out-Null

# --- Shared strings ------------------------------------------------------
# "Tabs" and "Missed Event" are brand-new strings; "Live Event", "Index",
# "0" and "1" already exist in the shared string table and get reused
# automatically.

Write-Host "List_Tab_Live_Event sheet created"
